$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

$ws.Range("A5").Value = "borxNoEquals"
$ws.Range("B5").Value = "BorxNo >= ,AND BorxNo <= "
$ws.Range("C5").Value = "CreateDate"

$ws.Range("A6").Value = "custNoAndBorxNo"
$ws.Range("B6").Value = "CustNo >= ,AND CustNo <= ,AND BorxNo >= ,AND BorxNo <= "

$ws.Activate()
$ws.Range("B7").Select()
